$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.898.56"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.259.00"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.84"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.256.94"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("E10").Value = "  +4.63%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.819.17"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.02"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.894.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E17").Value = "  +3.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.255.64"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.81%  "
$ws.Range("E22").Value = "  +5.75%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.76"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.510"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.394.47"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.19%  "
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.75"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.64"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.54"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.76"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.28"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.50"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.857"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.75%  "
$ws.Range("E40").Value = "  +9.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.01"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.748.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.43"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.88%  "
$ws.Range("E45").Value = "  +3.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "340.78"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.57%  "
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("E51").Value = "  +2.84%  "
